$d = $word.ActiveDocument

# Collapse "Anna Dimitrova " (the first two runs, before the bookmark) down to "1".
$d.Content.Find.Execute("Anna Dimitrova ", $true, $false, $false, $false, $false, $true, 1, $false, "1", 2)

# Remove the trailing "Nikolova" run that follows the bookmark.
$d.Content.Find.Execute("Nikolova", $true, $false, $false, $false, $false, $true, 1, $false, "", 2)
